$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("union-no-learning")
$ws.Activate()

# Row 4: C4 keeps its existing number-format style; D4:H4 get filled with real
# values and switch to a plain "font only" style (no custom number format).
$ws.Range("C4").Value = 0.59398496240601495

$rowFourPlain = $ws.Range("D4:H4")
$rowFourPlain.Style = "Normal"
$rowFourPlain.Font.Size = 13
$ws.Range("D4").Value = 0.0927318295739348
$ws.Range("E4").Value = 0.31328320802005
$ws.Range("F4").Value = 0.62656641604009999
$ws.Range("G4").Value = 0.21804511278195399
$ws.Range("H4").Value = 0.15538847117794399

# Row 8: all of C8:H8 switch to the same new plain style with real values.
$rowEight = $ws.Range("C8:H8")
$rowEight.Style = "Normal"
$rowEight.Font.Size = 13
$ws.Range("C8").Value = 0.360902255639097
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0.63909774436090205
$ws.Range("F8").Value = 0.39348370927318199
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0.60651629072681701

# Move the sheet's active selection to reflect where the author left off.
$ws.Range("L14").Select()
